# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt" /
# Cebollín at row 439, shifting the existing rows 439-476 down to 440-477
# (the former row 476 ends up as the new row 477).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 439 down by one row.
$ws.Rows.Item(439).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(439, 1).Value = 4
$ws.Cells.Item(439, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(439, 3).Value = "Los Lagos"
$ws.Cells.Item(439, 4).Value = 45106
$ws.Cells.Item(439, 5).Value = 10
$ws.Cells.Item(439, 6).Value = 100112037
$ws.Cells.Item(439, 7).Value = "Cebollín"
$ws.Cells.Item(439, 8).Value = "Sin especificar"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 60
$ws.Cells.Item(439, 11).Value = 7000
$ws.Cells.Item(439, 12).Value = 7000
$ws.Cells.Item(439, 13).Value = 7000
$ws.Cells.Item(439, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(439, 15).Value = "Región Metropolitana"
$ws.Cells.Item(439, 16).Value = 194
$ws.Cells.Item(439, 17).Value = 36
$ws.Cells.Item(439, 18).Value = "Hortaliza"

# Note: Rows.Item(439).Insert() already copies row 438's formatting down
# onto the newly created row 439, so column D keeps its date number
# format (style index 2) without any extra styling step here.
